$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.875.38"
$ws.Range("E2").Value = "  +2.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.904.02"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.62"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5009"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3001"
$ws.Range("E8").Value = "  +2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06863"
$ws.Range("E9").Value = "  +3.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.907.08"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.36"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07355"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.78"
$ws.Range("E13").Value = "  +7.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.131"
$ws.Range("E14").Value = "  +6.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6826"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.864.70"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008059"
$ws.Range("E17").Value = "  +2.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.37"
$ws.Range("E18").Value = "  +4.72%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.150.00"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.885"
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "184.51"
$ws.Range("E23").Value = "  +36.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.118"
$ws.Range("E24").Value = "  +9.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.397"
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.29"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.62"
$ws.Range("E27").Value = "  +11.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.956"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.398"
$ws.Range("E30").Value = "  +5.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08998"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.083"
$ws.Range("E32").Value = "  +3.70%  "
$ws.Range("E33").Value = "  +5.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7469"
$ws.Range("E34").Value = "  +5.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.145"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.669"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01930"
$ws.Range("E37").Value = "  +17.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.730"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.195"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9432"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4406"
$ws.Range("E41").Value = "  +5.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.34"
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.866"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.805"
$ws.Range("E45").Value = "  +4.29%  "
$ws.Range("E46").Value = "  +8.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05859"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3938"
$ws.Range("E48").Value = "  +6.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.596"
$ws.Range("E49").Value = "  +3.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.49"
$ws.Range("E50").Value = "  +3.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.400"
$ws.Range("E51").Value = "  +4.31%  "
